$d = $word.ActiveDocument

# Korrektur Bieber -> Biber
$d.Content.Find.Execute("Bea Bieber", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Bea Biber", 2)
